# Applies the "Optuna Attempt (go back with original)" edit to the workbook.
# Updates forecast values on the "Forecast Comparison" sheet and the
# corresponding rollup values on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------------
# Row -> hashtable of column letter -> new value
$forecastUpdates = @{
    2  = @{ D = 5;  H = 13;                 L = 1.04 }
    3  = @{ D = 5;  H = 10.71;              L = 0.85 }
    4  = @{ D = 5;  H = 9.710000000000001;  L = 0.95 }
    5  = @{ D = 5;  H = 8.710000000000001;  L = 1.06 }
    6  = @{ D = 10; H = 4.15;               L = 0.8100000000000001 }
    7  = @{         H = 5.47;               L = 1.14 }
    8  = @{ D = 10; H = 2.58;               L = 1.2 }
    9  = @{ D = 8;  H = 2.05;               L = 0.92 }
    10 = @{ D = 4;  H = 1.71;               L = 1.01 }
    11 = @{ D = 7;  H = 0.5;  I = "High"; J = "Urgent"; L = 0.99 }
    12 = @{ D = 9;  H = 0;                  L = 1.18 }
    13 = @{ D = 9;                          L = 0.99 }
    14 = @{ D = 4;                          L = 0.82 }
    15 = @{ D = 4;                          L = 0.8100000000000001 }
    16 = @{ D = 4;                          L = 1.19 }
    17 = @{ D = 7;                          L = 1.16 }
}

foreach ($row in $forecastUpdates.Keys) {
    $cols = $forecastUpdates[$row]
    foreach ($col in $cols.Keys) {
        $wsForecast.Range("$col$row").Value = $cols[$col]
    }
}

# --- Summary sheet --------------------------------------------------------------
$summaryUpdates = @{
    9  = "108"
    10 = "57"
    11 = "22"
    12 = "10"
    14 = "4"
}

foreach ($row in $summaryUpdates.Keys) {
    $cell = $wsSummary.Range("B$row")
    # Force the numeric-looking value to be stored as text, matching the
    # source workbook where these summary values are inline strings.
    $cell.NumberFormat = "@"
    $cell.Value = $summaryUpdates[$row]
}
